# BOT; UPDATE DATA
# Adds the newest day's (2020-05-21, serial 43972) infection figures to the
# three data sheets ("all", "kobe", "other"), refreshes the footnote that
# lists the out-of-city cases, and updates which sheet/tab is active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": append new row for 43972, pushing the footnote row down
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows(44).Insert()
$wsAll.Range("A44").Value = 43972
$wsAll.Range("B44").Value = 285
$wsAll.Range("C44").Value = 282
$wsAll.Range("D44").Value = 40
$wsAll.Range("E44").Value = 35
$wsAll.Range("F44").Value = 5
$wsAll.Range("G44").Value = 11
$wsAll.Range("H44").Value = 231

# ---------------------------------------------------------------------
# Sheet "kobe": correct the previous day's totals, then append new row
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Range("D98").Value = 2
$wsKobe.Range("E98").Value = 285

$wsKobe.Rows(99).Insert()
$wsKobe.Range("A99").Value = 43972
$wsKobe.Range("B99").Value = 0
$wsKobe.Range("C99").Value = 2956
$wsKobe.Range("D99").Value = 0
$wsKobe.Range("E99").Value = 285
$wsKobe.Range("F99").Value = 35
$wsKobe.Range("G99").Value = 31
$wsKobe.Range("H99").Value = 4
$wsKobe.Range("I99").Value = 11
$wsKobe.Range("J99").Value = 222

# ---------------------------------------------------------------------
# Sheet "other": totals for outside-Kobe cases are unchanged day over
# day, but a new row is still appended to keep the daily series going
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows(74).Insert()
$wsOther.Range("A74").Value = 43972
$wsOther.Range("B74").Value = 0
$wsOther.Range("C74").Value = 14
$wsOther.Range("D74").Value = 5
$wsOther.Range("E74").Value = 4
$wsOther.Range("F74").Value = 1
$wsOther.Range("G74").Value = 0
$wsOther.Range("H74").Value = 9

# ---------------------------------------------------------------------
# Footnote on "all"/"kobe" listing out-of-city patients gained one more
# case (no. 285), so the running count goes from 16 to 17
# ---------------------------------------------------------------------
$newNote = "※　24・34・53・58・59・60・158・161・163・192・237・248・268・272・276・277・285例目（計17件）は市外在住者です。"
$wsAll.Range("B45").Value = $newNote
$wsKobe.Range("B100").Value = $newNote

# ---------------------------------------------------------------------
# View state: make "all" the active sheet/tab, and update each sheet's
# active selection to sit on the freshly added rows
# ---------------------------------------------------------------------
$wsAll.Activate() | Out-Null
$wsAll.Range("C45").Select() | Out-Null

$wsKobe.Activate() | Out-Null
$unionSel = $excel.Union($wsKobe.Range("I99:J99"), $wsKobe.Range("F99"))
$unionSel.Select() | Out-Null

$wsOther.Activate() | Out-Null
$wsOther.Range("A73").Select() | Out-Null

$wsAll.Activate() | Out-Null
